# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted at row 32 (pushing the existing
# rows 32-40 down to 33-41). The new row carries the same market/category
# metadata as the rest of the sheet, with its own date and volume figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 32, shifting rows 32:40 down
# to 33:41 (dimension grows from A1:R40 to A1:R41 automatically).
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32.
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C32").Value = 'Ñuble'
$ws.Range("D32").Value = 44784
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112001
$ws.Range("G32").Value = 'Berenjena'
$ws.Range("H32").Value = 'Sin especificar'
$ws.Range("I32").Value = 'Primera'
$ws.Range("J32").Value = 80
$ws.Range("K32").Value = 12000
$ws.Range("L32").Value = 13000
$ws.Range("M32").Value = 12500
$ws.Range("N32").Value = '$/caja 60 unidades'
$ws.Range("O32").Value = 'Región de Arica y Parinacota'
$ws.Range("P32").Value = 208
$ws.Range("Q32").Value = 60
$ws.Range("R32").Value = 'Hortaliza'
